$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Content changes: the report is filled in as "admitted" (ADMITIDO) with a
#    real applicant's data instead of the placeholder "No Admitidos" example.
# ---------------------------------------------------------------------------

# Title row (row 6): "Listado No Admitidos" -> "Listado de Admitidos"
$ws.Range("B6").Value = "Listado de Admitidos"

# Data row (row 8): Nombre / Apellidos / Documento / SNP / Estado
$ws.Range("C8").Value = "DIEGO ALEXANDER"
$ws.Range("D8").Value = "CASTELLANOS JIMENEZ"

# Documento (E8) is an all-digit code; force it to stay text (as it already
# was) instead of being auto-converted to a number.
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "80762016"
$ws.Range("E8").NumberFormat = "General"

$ws.Range("F8").Value = "2023MCHB2077"
$ws.Range("J8").Value = "ADMITIDO"

# ---------------------------------------------------------------------------
# 2) Formatting changes: the whole header/table area is switched to a Text
#    ("@") number format (so codes/ids typed later keep their exact digits),
#    while keeping the existing borders / fills / alignment untouched.
# ---------------------------------------------------------------------------

$ws.Range("A1:A5").NumberFormat = "@"

$ws.Range("B1:J5").NumberFormat = "@"
$ws.Range("B1:J5").WrapText = $true

$ws.Range("A6").NumberFormat = "@"

$ws.Range("B6:I6").NumberFormat = "@"

$ws.Range("J6").NumberFormat = "@"

$ws.Range("A7:J7").NumberFormat = "@"

# ---------------------------------------------------------------------------
# 3) Column sizing: columns A-F and J get the sheet's default width, matching
#    the new layout that frames the table (G:I keep their existing widths).
# ---------------------------------------------------------------------------

$ws.Columns("A:F").ColumnWidth = 10.666666666666666
$ws.Columns("J").ColumnWidth = 10.666666666666666

# ---------------------------------------------------------------------------
# 4) Selection: land on the whole sheet (A1:J1048576) instead of a stray cell.
# ---------------------------------------------------------------------------

$ws.Range("A1:J1048576").Select()
